$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 35

$ws.Cells.Item($row, 1).Value = 45881
$ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item($row, 2).Value = 6441.64013671875
$ws.Cells.Item($row, 3).Value = 6444.25
$ws.Cells.Item($row, 4).Value = 6385.759765625
$ws.Cells.Item($row, 5).Value = 6395.169921875
$ws.Cells.Item($row, 6).Value = 1667322000
$ws.Cells.Item($row, 7).Value = 0.0106990624099332
$ws.Cells.Item($row, 8).Value = 1

$textI = @"
from power plant to the moon the u nuclear revival in dogecoin jump in a week what is going on teenage mutant ninja turtle pizzeria to debut in la thestreet pro analyst tackle tariff crypto and more in market roundtable solar stock outshine sp in despite policy headwind why share of tesla are soaring today time tracking and trust what today employee really think amazon is selling a storage cart for thats great for small place cold wallet m presale entry beat xrps price test and eths resistance goal why many smallbusinesses owner fear even a short recession will sink them creating a load history tracker that improves rate negotiation crypto news today strategy m bitcoin purchase bullish upsizes ipo alt sigma b fundraising move solana price prediction etf approval could help sol overtake xrp but this cheaper token will rally x and outshine their roi bank of america strategist expect fed to lower interest rate benny brown acquires wide open music publishing form new publishing company affirm stock revisits this place a buy now pay later delight shopper demi lovato reunites with exboyfriend joe jonas for a camp rock performance new pickleball kingdom indoor facility coming soon to south brunswick nj nasfm honor fire chief todd parker for excellence in emergency response virginia obamacare premium set to jump sharply scc report virginia obamacare premium set to jump sharply scc report going online in russia can be frustrating complicated and even dangerous chatgpt get smarter with gpt top feature youll love if youd invested in interactive broker ibkr stock year ago here how much youd have today graphene market to grow at cagr by lynch carpenter investigates claim in pandora data breach neptune frozen treat opening in malaga cove plaza top tech news today taylor swift ai scandal iphone price hike tesla dojo shutdown metas talent war siri upgrade stock market party to be ruined by sudden economic slowdown say stifel cnbc bitcoin hovers around whats next lockheed martin could lose major deal after u hit switzerland with tariff top meme coin to invest for maximum gain in xrp news today analyst react to m xrp transfer to coinbase is it routine move ford hit the pedal on ev production with billion overhaul of kentucky plant fed bowman say job data solidifies case for rate cut u consumer to bear brunt of tariff hit goldman economist say u consumer to bear brunt of tariff hit goldman economist say roll royce share price is on a sizzling hot streak but beware of the danger top laptop brand in the usa in ibd live qa stock list for monday aug why good leader ask great question manufacturer urged to act with technologyfirst plan a disruption intensifies say infotech research group in new report philippe cousteau earthecho international announces diverse youth advocate join the blue carbon ambassador program play the perfect game this summer with flightpath golf tee what to expect from tuesday report on inflation investopedia evercommerce to present at third quarter investor conference evercommerce to present at third quarter investor conference comedic multihyphenate matt mathews sign with wme hurry and grab a macbook air for just at back market air india suspends flight to washington dc amid ongoing tension this is what whale are betting on lockheed martin smart money is betting big in ups option aircraft gasket market is forecasted to reach u million in say stratview research best cryptos to invest in that could deliver massive longterm gain blockdag litecoin stellar toncoin solar stock outshine oil and gas benchmark after climate credit cut aerospace pilot control market is forecasted to reach u billion in say stratview research xml announces launch of new xml woman website relaunch of xml woman division aircraft gasket market is forecasted to reach u million in say stratview research aerospace pilot control market is forecasted to reach u billion in say stratview research diversified foodservice supply appoints lev peker a chief executive officer body armor plate market is forecasted to reach u million in say stratview research body armor plate market is forecasted to reach u million in say stratview research best typec hub for macbook in venchi bringing beloved italian chocolate to century city tegna stock rocket on nexstar merger talk fcc appears set to ease local tv rule in order to smooth deal path deadline btc miner stance on cloud mining in if youd invested in zoetis zts stock year ago here how much youd have today fda may not renew pfizers authorization for covid shot in younger child moderna prepares to boost supply ja mining offer fixed usd yield with crypto deposit alphabet inc nasdaqgoog share acquired by gouws capital llc kohmann bosshard financial service llc purchase share of costco wholesale corporation nasdaqcost gouws capital llc boost stock holding in thermo fisher scientific inc nysetmo ana alsharif became the face of the war in gaza for million then israel killed him home condition made rescue effort difficult pfd protect partner with sentin to pioneer aipowered industrial inspection gear up for the new school year with off therabody best seller wall street biggest bear is expecting a stockmarket pullback this year trump say gold import will not face u tariff bloombergcom irobots future isnt looking up russian seeking asylum in u facing new obstacle say usbased activist morse corp awarded m u army ota for uas development morse corp awarded m u army ota for uas development aem united state ag tractor and combine report july aem united state ag tractor and combine report july scientist say they cruised the ocean in a deepsea submersible and came across an undiscovered ecosystem macys is selling a cozy puma hoodie for thats super soft winnebago industry announces strategic leadership change to drive future growth how to stretch backtoschool budget at salvage store openai and google vet periodic lab seek million investment aol dialup got america online in the s it ending next month trump say gold will not be tariffed cnbc how turmeric water fight inflammation aussivo reveals vision to build blockchainverified cloud trust layer a new standard for transparency in the digital world the viral govee floor lamp is still at prime day price silver at doublevibration support ahead of gann time cluster unit building at rockaway reach completion in brownsville fill up your shopping cart with today best target deal including over off outdoor gear headphone kitchen good and more communication service stock whale activity in today session playhouse under construction for popup park in virginia beach aol is finally disconnecting it dialup internet year on barack obamas official portrait moved to nonpublic area in white house western union ceo ice raid tariff no threat to instamex deal video show moment of explosion at u steel plant in pennsylvania trump gross white house makeover reach petty new low crowded north county field vying to replace brian williams in missouri senate how the world most valuable company got caught in the middle of trump spat with china commentary vance right to call out warped partisan representation why is trump attacking intels ceo demi lovato reunites with exboyfriend joe jonas for a camp rock performance mamdani launch tour of new york city with a message linking cuomo to trump nasfm honor fire chief todd parker for excellence in emergency response what to know about the putintrump summit in alaska what to know about the putintrump summit in alaska what to know about the putintrump summit in alaska what to know about the putintrump summit in alaska what to know about the putintrump summit in alaska what to know about the putintrump summit in alaska what to know about the putintrump summit in alaska what to know about the putintrump summit in alaska what to know about the putintrump summit in alaska what to know about the putintrump summit in alaska what to know about the putintrump summit in alaska trump lash out at crockett renews call for cognitive test lynch carpenter investigates claim in pandora data breach president mayor police chief social worker trump is a man of many hat trump say he and putin will discus ukraine future border israel plan to widen coming offensive beyond gaza city into last area not under it control israel plan to widen coming offensive beyond gaza city into last area not under it control israel plan to widen coming offensive beyond gaza city into last area not under it control u senate candidate set the tone for kentucky republican during statewide dinner what to know about the putintrump summit in alaska what to know about the putintrump summit in alaska photo of trump emergency declaration a the national guard is activated in washington photo of trump emergency declaration a the national guard is activated in washington trump mistake alaska for russia and yes there snark surveillance case test chapter a standing and unfairness standard clinton county sheriff office arrest individual involved with multicounty sex trafficking operation a judge brutal rebuke of trump epstein gambit cnn politics what to know about the putintrump summit in alaska local view give ailing minnesota kid a fighting chance philippe cousteau earthecho international announces diverse youth advocate join the blue carbon ambassador program california couple file lawsuit against klp in conway over shortterm rental rule rubio hit medium over gaza narrative say some ally making it harder to achieve peace multiple men have impersonated ice agent to kidnap and assault woman multiple men have impersonated ice agent to kidnap and assault woman doe using ai dumb you down vote to save a historic building from demolition trump promised to fix inflation he doing this instead south carolina lawmaker awaiting trial on child sex abuse material charge resigns from office nvidia to share of revenue from chip sale to china with u government trump say gold wont face tariff trump say gold wont face tariff trump say gold wont face tariff trump say gold wont face tariff trump say gold wont face tariff ana alsharif became the face of the war in gaza for million then israel killed him trump say he placing washington police under federal control and deploying the national guard trump declares crime in dc will be cleaned up real quick during his announcement that he taking over the city police department chinese ship collide on chase in disputed water ap news summary at pm edt trump say he placing washington police under federal control and activating the national guard houlihan financial resource group ltd sell share of international business machine corporation nyseibm meridian wealth management llc ha million position in broadcom inc nasdaqavgo motley fool asset management llc raise stock holding in broadcom inc nasdaqavgo colorado congressional democrat visit aurora ice facility meet with detained activist jeanette vizguerra colorado congressional democrat visit aurora ice facility meet with detained activist jeanette vizguerra colorado congressional democrat visit aurora ice facility meet with detained activist jeanette vizguerra colorado congressional democrat visit aurora ice facility meet with detained activist jeanette vizguerra russian seeking asylum in u facing new obstacle say usbased activist trump warns protester amid dc crime crack down you spit and we hit trump suggests he can send national guard troop into other u city beyond dc scientist say they cruised the ocean in a deepsea submersible and came across an undiscovered ecosystem dillon associate inc increase holding in salesforce inc nysecrm colombian presidential candidate miguel uribe turbay dy two month after assassination attempt live now texas hold legislative session after democrat leave state jeffries hammer trump for seizing control of dc police no basis jeffries hammer trump for seizing control of dc police no basis jeffries hammer trump for seizing control of dc police no basis jeffries hammer trump for seizing control of dc police no basis maga v miga laura loomer block marjorie taylor greene on social medium escalating a feud over america israel policy trump looking at reclassifying marijuana a less dangerous drug trump looking at reclassifying marijuana a less dangerous drug danish zoo asks for donated pet to feed to carnivore trump looking at reclassifying marijuana a less dangerous drug trump looking at reclassifying marijuana a less dangerous drug trump looking at reclassifying marijuana a less dangerous drug el paso police man dy when he hit tree get ejected trump team get blasted by judge for clumsy epstein diversion raw politics is now driving law enforcement at the department of justice civil right division nancy pelosi get all high and mighty lying about trump j and the ntl guard and oh hell no old woman this candidate for california governor ha a potential conflict of interest in her own home ishares core sp etf nysearcaivv share sold by aveo capital partner llc maine state library not victimized a suggested letter maine state library not victimized a suggested letter here are the country that have recognized a palestinian state johnson fistel investigates fairness of proposed sale of meridianlink currently the most dangerous place maia wealth llc purchase share of ishares core u aggregate bond etf nysearcaagg my wife of year just died my grief is so overwhelming i can barely cope
"@

$textJ = @"
power plant moon nuclear revival dogecoin jump teenage mutant ninja turtle pizzeria debut la thestreet pro analyst tackle tariff crypto market roundtable solar stock outshine sp policy headwind tesla soaring tracking trust employee amazon selling storage cart thats great small place cold wallet presale entry beat xrps price test eths resistance goal smallbusinesses owner fear short recession sink creating load history tracker improves rate negotiation crypto strategy bitcoin purchase bullish upsizes ipo alt sigma b fundraising move solana price prediction approval help sol overtake xrp cheaper token rally x outshine roi bank america strategist expect fed lower interest rate benny brown acquires wide open music publishing form publishing company affirm stock revisits place buy pay later delight shopper demi lovato reunites exboyfriend joe jonas camp rock performance pickleball kingdom indoor facility coming soon south brunswick nj nasfm honor fire chief todd parker excellence emergency response virginia obamacare premium set jump sharply scc report virginia obamacare premium set jump sharply scc report online russia frustrating complicated dangerous chatgpt smarter gpt feature youll love youd invested interactive broker ibkr stock ago youd graphene market grow cagr lynch carpenter investigates claim pandora data breach neptune frozen treat opening malaga cove plaza tech taylor swift ai scandal iphone price hike tesla dojo shutdown metas talent war siri upgrade stock market party ruined sudden economic slowdown stifel cnbc bitcoin hovers around whats next lockheed martin lose major deal hit switzerland tariff meme coin invest maximum gain xrp analyst react xrp transfer coinbase routine move ford hit pedal ev production billion overhaul kentucky plant fed bowman job data solidifies rate cut consumer bear brunt tariff hit goldman economist consumer bear brunt tariff hit goldman economist roll royce price sizzling hot streak beware danger laptop brand usa ibd qa stock list monday aug good great manufacturer urged act technologyfirst plan disruption intensifies infotech research report philippe cousteau earthecho international announces diverse youth advocate join blue carbon ambassador program play perfect game summer flightpath golf tee expect tuesday report inflation investopedia evercommerce present third quarter investor conference evercommerce present third quarter investor conference comedic multihyphenate matt mathews sign wme hurry grab macbook air market air india suspends flight washington dc amid ongoing tension whale betting lockheed martin smart money betting ups option aircraft gasket market forecasted reach million stratview research best cryptos invest deliver massive longterm gain blockdag litecoin stellar toncoin solar stock outshine oil gas benchmark climate credit cut aerospace pilot control market forecasted reach billion stratview research xml announces launch xml woman website relaunch xml woman division aircraft gasket market forecasted reach million stratview research aerospace pilot control market forecasted reach billion stratview research diversified foodservice supply appoints lev peker chief executive officer body armor plate market forecasted reach million stratview research body armor plate market forecasted reach million stratview research best typec hub macbook venchi bringing beloved italian chocolate century tegna stock rocket nexstar merger fcc set ease local tv rule order smooth deal path deadline btc miner stance cloud mining youd invested zoetis zts stock ago youd fda renew pfizers authorization covid shot younger child moderna prepares boost supply ja mining offer fixed usd yield crypto deposit alphabet nasdaqgoog acquired gouws kohmann bosshard financial service purchase costco wholesale corporation nasdaqcost gouws boost stock thermo fisher scientific nysetmo ana alsharif became face war gaza million israel killed condition rescue effort difficult pfd protect partner sentin pioneer aipowered industrial inspection gear school therabody best seller wall street biggest bear expecting stockmarket pullback trump gold import face tariff bloombergcom irobots future isnt russian seeking asylum facing obstacle usbased activist morse corp awarded army ota uas development morse corp awarded army ota uas development aem united ag tractor combine report aem united ag tractor combine report scientist cruised ocean deepsea submersible came across undiscovered ecosystem macys selling cozy puma hoodie thats super soft winnebago industry announces strategic leadership change drive future growth stretch backtoschool budget salvage store openai google vet periodic lab seek million investment aol dialup got america online ending next month trump gold tariffed cnbc turmeric water fight inflammation aussivo reveals vision build blockchainverified cloud trust layer transparency digital world viral govee floor lamp prime price silver doublevibration support ahead gann cluster unit building rockaway reach completion brownsville fill shopping cart best target deal outdoor gear headphone kitchen good communication service stock whale activity session playhouse construction popup park virginia beach aol finally disconnecting dialup internet barack obamas portrait moved nonpublic area white house western union ceo raid tariff threat instamex deal video moment explosion steel plant pennsylvania trump gross white house makeover reach petty low crowded north county field vying replace brian williams missouri senate world valuable company got caught middle trump spat china commentary vance call warped partisan representation trump attacking intels ceo demi lovato reunites exboyfriend joe jonas camp rock performance mamdani launch tour york message linking cuomo trump nasfm honor fire chief todd parker excellence emergency response putintrump summit alaska putintrump summit alaska putintrump summit alaska putintrump summit alaska putintrump summit alaska putintrump summit alaska putintrump summit alaska putintrump summit alaska putintrump summit alaska putintrump summit alaska putintrump summit alaska trump lash crockett renews call cognitive test lynch carpenter investigates claim pandora data breach president mayor police chief social worker trump man hat trump putin discus ukraine future border israel plan widen coming offensive beyond gaza area control israel plan widen coming offensive beyond gaza area control israel plan widen coming offensive beyond gaza area control senate candidate set tone kentucky republican statewide dinner putintrump summit alaska putintrump summit alaska photo trump emergency declaration national guard activated washington photo trump emergency declaration national guard activated washington trump mistake alaska russia yes snark surveillance test chapter standing unfairness clinton county sheriff office arrest individual involved multicounty sex trafficking operation judge brutal rebuke trump epstein gambit cnn politics putintrump summit alaska local view give ailing minnesota kid fighting chance philippe cousteau earthecho international announces diverse youth advocate join blue carbon ambassador program california couple file lawsuit klp conway shortterm rental rule rubio hit medium gaza narrative ally making harder achieve peace multiple men impersonated agent kidnap assault woman multiple men impersonated agent kidnap assault woman doe ai dumb vote save historic building demolition trump promised fix inflation instead south carolina lawmaker awaiting trial child sex abuse material charge resigns office nvidia revenue chip sale china government trump gold wont face tariff trump gold wont face tariff trump gold wont face tariff trump gold wont face tariff trump gold wont face tariff ana alsharif became face war gaza million israel killed trump placing washington police federal control deploying national guard trump declares crime dc cleaned real quick announcement taking police department chinese ship collide chase disputed water summary pm trump placing washington police federal control activating national guard houlihan financial resource ltd sell international business machine corporation nyseibm meridian wealth million broadcom nasdaqavgo motley fool asset raise stock broadcom nasdaqavgo colorado congressional democrat visit aurora facility meet detained activist jeanette vizguerra colorado congressional democrat visit aurora facility meet detained activist jeanette vizguerra colorado congressional democrat visit aurora facility meet detained activist jeanette vizguerra colorado congressional democrat visit aurora facility meet detained activist jeanette vizguerra russian seeking asylum facing obstacle usbased activist trump warns protester amid dc crime crack spit hit trump suggests send national guard troop beyond dc scientist cruised ocean deepsea submersible came across undiscovered ecosystem dillon associate increase salesforce nysecrm colombian presidential candidate miguel uribe turbay dy month assassination attempt texas legislative session democrat jeffries hammer trump seizing control dc police basis jeffries hammer trump seizing control dc police basis jeffries hammer trump seizing control dc police basis jeffries hammer trump seizing control dc police basis jeffries hammer trump seizing control dc police basis maga v miga laura loomer block marjorie taylor greene social medium escalating feud america israel policy trump reclassifying marijuana less dangerous drug trump reclassifying marijuana less dangerous drug danish zoo asks donated pet feed carnivore trump reclassifying marijuana less dangerous drug trump reclassifying marijuana less dangerous drug trump reclassifying marijuana less dangerous drug el paso police man dy hit tree ejected trump team blasted judge clumsy epstein diversion raw politics driving enforcement department justice civil division nancy pelosi high mighty lying trump j ntl guard oh hell old woman candidate california governor potential conflict interest ishares core sp nysearcaivv aveo partner maine library victimized suggested maine library victimized suggested country recognized palestinian johnson fistel investigates fairness proposed sale meridianlink currently dangerous place maia wealth purchase ishares core aggregate bond nysearcaagg wife died grief overwhelming barely cope
"@

$ws.Cells.Item($row, 9).Value = $textI
$ws.Cells.Item($row, 10).Value = $textJ
